$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a literal text value into a cell without Excel's automatic
# number/date conversion (e.g. "0.00" or "2021-03-17"), while preserving the
# cell's existing style (so no new style entries get created).
# We do this by computing the literal text via a TEXT() formula in a scratch
# cell, copying the destination's own style back onto itself, then pasting
# just the computed text value (not the formula) into the destination.
function Set-LiteralText {
    param($Address, $Text)

    $scratch = $ws.Range("ZZ1")
    $escapedText = $Text.Replace('"', '""')
    $scratch.Formula = '=TEXT("' + $escapedText + '","")'

    $dest = $ws.Range($Address)

    # Preserve destination's current formatting.
    $dest.Copy() | Out-Null
    $dest.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    # Bring over the literal text as a value only (keeps style untouched).
    $scratch.Copy() | Out-Null
    $dest.PasteSpecial(-4163) | Out-Null   # xlPasteValues

    $scratch.ClearContents() | Out-Null
}

# Row 1 (headers): BillingAddress is inserted, ZipCode keeps its place right
# after it, and StreetAddress moves to the new last column (E).
$ws.Range("C1").Value = "BillingAddress"
$ws.Range("D1").Value = "ZipCode"
$ws.Range("E1").Value = "StreetAddress"

# Row 2: clear the old values that used to live here (they move to row 3),
# leaving blank placeholders under InvoiceDate/BillingAddress, and record the
# extracted zip code value under the ZipCode column.
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
Set-LiteralText "D2" "0.00"

# Row 3: the invoice date and the full billing address text that used to sit
# in row 2 now live in row 3, under InvoiceDate (B) and StreetAddress (E).
Set-LiteralText "B3" "2021-03-17"
$ws.Range("E3").Value = "City, State, Country: [City, State, Country]"
